$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 475: remove B475:D475 (previously empty inline-string placeholder cells) ---
$ws.Range("B475:D475").ClearContents()

# --- Row 476 ---
$ws.Range("A476").Value = "Hooklight 2"
$ws.Range("F476").NumberFormat = "@"
$ws.Range("F476").Value = "2"
$ws.Range("G476").NumberFormat = "@"
$ws.Range("G476").Value = "12"
$ws.Range("H476").Value = 229
$ws.Range("I476").Value = "5th Dr"
$ws.Range("J476").Value = "Adric, Nyssa, Tegan"
$ws.Range("K476").Value = "8th Dr"
$ws.Range("L476").Value = "Nura, Nigh Guard, Halcyon"
$ws.Range("M476").Value = "Tim Foley"
$ws.Range("N476").Value = "Ken Bentley"
$ws.Range("O476").NumberFormat = "@"
$ws.Range("O476").Value = "2025"

# --- Row 477 ---
$ws.Range("A477").Value = "The Story Demon"
$ws.Range("E477").Value = "The Cosmos and Mrs Clarke"
$ws.Range("F477").NumberFormat = "@"
$ws.Range("F477").Value = "1"
$ws.Range("G477").NumberFormat = "@"
$ws.Range("G477").Value = "2"
$ws.Range("H477").Value = 82
$ws.Range("I477").Value = "6th Dr"
$ws.Range("J477").Value = "Constance"
$ws.Range("L477").Value = "Dalek mutant"
$ws.Range("M477").Value = "Julian Richards"
$ws.Range("N477").Value = "Samuel Clemens"
$ws.Range("O477").NumberFormat = "@"
$ws.Range("O477").Value = "2025"

# --- Row 478 ---
$ws.Range("A478").Value = "The Key to Many Worlds"
$ws.Range("E478").Value = "The Cosmos and Mrs Clarke"
$ws.Range("F478").NumberFormat = "@"
$ws.Range("F478").Value = "2"
$ws.Range("G478").NumberFormat = "@"
$ws.Range("G478").Value = "2"
$ws.Range("H478").Value = 88
$ws.Range("I478").Value = "6th Dr"
$ws.Range("J478").Value = "Constance"
$ws.Range("K478").Value = "Iris, Marco"
$ws.Range("M478").Value = "Paul Magrs"
$ws.Range("N478").Value = "Samuel Clemens"
$ws.Range("O478").NumberFormat = "@"
$ws.Range("O478").Value = "2025"

# --- Row 479 ---
$ws.Range("A479").Value = "Inconstancy"
$ws.Range("E479").Value = "The Cosmos and Mrs Clarke"
$ws.Range("F479").NumberFormat = "@"
$ws.Range("F479").Value = "3"
$ws.Range("G479").NumberFormat = "@"
$ws.Range("G479").Value = "2"
$ws.Range("H479").Value = 82
$ws.Range("I479").Value = "6th Dr"
$ws.Range("J479").Value = "Constance"
$ws.Range("L479").Value = "Claudia Purnell"
$ws.Range("M479").Value = "Ian Potter"
$ws.Range("N479").Value = "Samuel Clemens"
$ws.Range("O479").NumberFormat = "@"
$ws.Range("O479").Value = "2025"

# --- Row 480 ---
$ws.Range("A480").Value = "Missy Part 2"
$ws.Range("E480").Value = "N/A"
$ws.Range("F480").NumberFormat = "@"
$ws.Range("F480").Value = "2"
$ws.Range("G480").Value = "N/A"
$ws.Range("H480").Value = 77
$ws.Range("I480").Value = "N, /, A"
$ws.Range("J480").Value = "N/A"
$ws.Range("K480").Value = "N, /, A"
$ws.Range("L480").Value = "N/A"
$ws.Range("M480").Value = "N/A"
$ws.Range("N480").Value = "N/A"
$ws.Range("O480").NumberFormat = "@"
$ws.Range("O480").Value = "2025"

# --- Row 481 ---
$ws.Range("A481").Value = "A Forest of All Seasons"
$ws.Range("E481").Value = "A Feast of Steven"
$ws.Range("F481").NumberFormat = "@"
$ws.Range("F481").Value = "3"
$ws.Range("G481").NumberFormat = "@"
$ws.Range("G481").Value = "1"
$ws.Range("H481").Value = 45
$ws.Range("I481").Value = "1st Dr"
$ws.Range("J481").Value = "Steven, Vicki"
$ws.Range("M481").Value = "Jacqueline Rayner"
$ws.Range("N481").Value = "John Ainsworth"
$ws.Range("O481").NumberFormat = "@"
$ws.Range("O481").Value = "2025"

# --- Row 482 ---
$ws.Range("A482").Value = "The Doctor's Gambit"
$ws.Range("E482").Value = "A Feast of Steven"
$ws.Range("F482").NumberFormat = "@"
$ws.Range("F482").Value = "4"
$ws.Range("G482").NumberFormat = "@"
$ws.Range("G482").Value = "1"
$ws.Range("H482").Value = 44
$ws.Range("I482").Value = "Steven, Dodo"
$ws.Range("K482").Value = "1st Dr"
$ws.Range("M482").Value = "Jacqueline Rayner"
$ws.Range("N482").Value = "John Ainsworth"
$ws.Range("O482").NumberFormat = "@"
$ws.Range("O482").Value = "2025"

# --- Row 483 ---
$ws.Range("A483").Value = "A Feast of Steven - Behind the Scenes"
$ws.Range("E483").Value = "N/A"
$ws.Range("F483").NumberFormat = "@"
$ws.Range("F483").Value = "5"
$ws.Range("G483").Value = "N/A"
$ws.Range("H483").Value = 26
$ws.Range("I483").Value = "N, /, A"
$ws.Range("J483").Value = "N/A"
$ws.Range("K483").Value = "N, /, A"
$ws.Range("L483").Value = "N/A"
$ws.Range("M483").Value = "N/A"
$ws.Range("N483").Value = "N/A"
$ws.Range("O483").NumberFormat = "@"
$ws.Range("O483").Value = "2025"

# --- Row 484 ---
$ws.Range("A484").Value = "The Remains of Kaerula"
$ws.Range("E484").Value = "The Ruins of Kaerula"
$ws.Range("F484").NumberFormat = "@"
$ws.Range("F484").Value = "3"
$ws.Range("H484").Value = 77
$ws.Range("I484").Value = "4th Dr"
$ws.Range("J484").Value = "Leela, K9"
$ws.Range("M484").Value = "Phil Mulryne"
$ws.Range("N484").Value = "Helen Goldwyn, Nicholas Briggs, Jamie Anderson"
$ws.Range("O484").NumberFormat = "@"
$ws.Range("O484").Value = "2025"

# --- Row 485 ---
$ws.Range("A485").Value = "The Ruins of Kaerula"
$ws.Range("E485").Value = "The Ruins of Kaerula"
$ws.Range("F485").NumberFormat = "@"
$ws.Range("F485").Value = "4"
$ws.Range("H485").Value = 77
$ws.Range("I485").Value = "4th Dr"
$ws.Range("J485").Value = "Leela, K9"
$ws.Range("M485").Value = "Phil Mulryne"
$ws.Range("N485").Value = "Helen Goldwyn, Nicholas Briggs, Jamie Anderson"
$ws.Range("O485").NumberFormat = "@"
$ws.Range("O485").Value = "2025"

# --- Row 486 ---
$ws.Range("A486").Value = "Cry of the Banshee"
$ws.Range("E486").Value = "The Ruins of Kaerula"
$ws.Range("F486").NumberFormat = "@"
$ws.Range("F486").Value = "5"
$ws.Range("H486").Value = 75
$ws.Range("I486").Value = "4th Dr"
$ws.Range("J486").Value = "Leela, K9"
$ws.Range("M486").Value = "Tim Foley"
$ws.Range("N486").Value = "Helen Goldwyn, Nicholas Briggs, Jamie Anderson"
$ws.Range("O486").NumberFormat = "@"
$ws.Range("O486").Value = "2025"

# --- Row 487 ---
$ws.Range("A487").Value = "With the Angels Part 1-2"
$ws.Range("E487").Value = "N/A"
$ws.Range("F487").NumberFormat = "@"
$ws.Range("F487").Value = "1"
$ws.Range("G487").Value = "N/A"
$ws.Range("H487").Value = 64
$ws.Range("I487").Value = "N, /, A"
$ws.Range("J487").Value = "N/A"
$ws.Range("K487").Value = "N, /, A"
$ws.Range("L487").Value = "N/A"
$ws.Range("M487").Value = "N/A"
$ws.Range("N487").Value = "N/A"
$ws.Range("O487").NumberFormat = "@"
$ws.Range("O487").Value = "2025"

# --- Row 488 ---
$ws.Range("A488").Value = "Catastrophix"
$ws.Range("E488").Value = "Past Forward"
$ws.Range("F488").NumberFormat = "@"
$ws.Range("F488").Value = "2"
$ws.Range("G488").NumberFormat = "@"
$ws.Range("G488").Value = "2"
$ws.Range("H488").Value = 67
$ws.Range("I488").Value = "7th Dr"
$ws.Range("J488").Value = "Harry, Naomi, Ray"
$ws.Range("M488").Value = "Lizzie Hopley"
$ws.Range("N488").Value = "Samuel Clemens"
$ws.Range("O488").NumberFormat = "@"
$ws.Range("O488").Value = "2025"

# --- Row 489 ---
$ws.Range("A489").Value = "With the Angels Part 3-4"
$ws.Range("E489").Value = "N/A"
$ws.Range("F489").NumberFormat = "@"
$ws.Range("F489").Value = "3"
$ws.Range("G489").Value = "N/A"
$ws.Range("H489").Value = 67
$ws.Range("I489").Value = "N, /, A"
$ws.Range("J489").Value = "N/A"
$ws.Range("K489").Value = "N, /, A"
$ws.Range("L489").Value = "N/A"
$ws.Range("M489").Value = "N/A"
$ws.Range("N489").Value = "N/A"
$ws.Range("O489").NumberFormat = "@"
$ws.Range("O489").Value = "2025"

# --- Row 490 ---
$ws.Range("A490").Value = "Missy Part 3"
$ws.Range("E490").Value = "N/A"
$ws.Range("F490").NumberFormat = "@"
$ws.Range("F490").Value = "3"
$ws.Range("G490").Value = "N/A"
$ws.Range("H490").Value = 77
$ws.Range("I490").Value = "N, /, A"
$ws.Range("J490").Value = "N/A"
$ws.Range("K490").Value = "N, /, A"
$ws.Range("L490").Value = "N/A"
$ws.Range("M490").Value = "N/A"
$ws.Range("N490").Value = "N/A"
$ws.Range("O490").NumberFormat = "@"
$ws.Range("O490").Value = "2025"

# --- Row 491 ---
$ws.Range("A491").Value = "The Dead Sea"
$ws.Range("F491").NumberFormat = "@"
$ws.Range("F491").Value = "1"
$ws.Range("G491").NumberFormat = "@"
$ws.Range("G491").Value = "3"
$ws.Range("H491").Value = 64
$ws.Range("I491").Value = "War Dr"
$ws.Range("M491").Value = "Alfie Shaw"
$ws.Range("N491").Value = "Ken Bentley"
$ws.Range("O491").NumberFormat = "@"
$ws.Range("O491").Value = "2025"

# --- Row 492 ---
$ws.Range("A492").Value = "Unit 26"
$ws.Range("F492").NumberFormat = "@"
$ws.Range("F492").Value = "2"
$ws.Range("G492").NumberFormat = "@"
$ws.Range("G492").Value = "3"
$ws.Range("H492").Value = 69
$ws.Range("I492").Value = "War Dr"
$ws.Range("M492").Value = "Alfie Shaw"
$ws.Range("N492").Value = "Ken Bentley"
$ws.Range("O492").NumberFormat = "@"
$ws.Range("O492").Value = "2025"

# --- Row 493 ---
$ws.Range("A493").Value = "Yesterday is Tomorrow and Tomorrow is Today"
$ws.Range("F493").NumberFormat = "@"
$ws.Range("F493").Value = "3"
$ws.Range("G493").NumberFormat = "@"
$ws.Range("G493").Value = "3"
$ws.Range("H493").Value = 73
$ws.Range("I493").Value = "War Dr"
$ws.Range("M493").Value = "Alfie Shaw"
$ws.Range("N493").Value = "Ken Bentley"
$ws.Range("O493").NumberFormat = "@"
$ws.Range("O493").Value = "2025"

# --- Row 494 ---
$ws.Range("A494").Value = "Kaiju"
$ws.Range("E494").Value = "Fractures"
$ws.Range("F494").NumberFormat = "@"
$ws.Range("F494").Value = "1"
$ws.Range("H494").Value = 74
$ws.Range("I494").Value = "Bambera, Savarin, Rix"
$ws.Range("K494").Value = "McManis"
$ws.Range("M494").Value = "Robert Valentine"
$ws.Range("N494").Value = "Samuel Clemens"
$ws.Range("O494").NumberFormat = "@"
$ws.Range("O494").Value = "2025"

# --- Row 495 ---
$ws.Range("A495").Value = "Debrief"
$ws.Range("E495").Value = "Fractures"
$ws.Range("F495").NumberFormat = "@"
$ws.Range("F495").Value = "2"
$ws.Range("H495").Value = 75
$ws.Range("I495").Value = "Zbrigniev"
$ws.Range("K495").Value = "Bambera"
$ws.Range("L495").Value = "Brigade Leader, Winifred Bambera"
$ws.Range("M495").Value = "Alfie Shaw"
$ws.Range("N495").Value = "Samuel Clemens"
$ws.Range("O495").NumberFormat = "@"
$ws.Range("O495").Value = "2025"

# --- Row 496 ---
$ws.Range("A496").Value = "Shatterpoint"
$ws.Range("E496").Value = "Fractures"
$ws.Range("F496").NumberFormat = "@"
$ws.Range("F496").Value = "3"
$ws.Range("H496").Value = 72
$ws.Range("I496").Value = "Bambera, Savarin, Rix"
$ws.Range("M496").Value = "Mark Wright"
$ws.Range("N496").Value = "Samuel Clemens"
$ws.Range("O496").NumberFormat = "@"
$ws.Range("O496").Value = "2025"

# --- Row 497 ---
$ws.Range("A497").Value = "The Voord in London"
$ws.Range("E497").Value = "N/A"
$ws.Range("F497").NumberFormat = "@"
$ws.Range("F497").Value = "1"
$ws.Range("G497").Value = "N/A"
$ws.Range("H497").Value = 84
$ws.Range("I497").Value = "N, /, A"
$ws.Range("J497").Value = "N/A"
$ws.Range("K497").Value = "N, /, A"
$ws.Range("L497").Value = "N/A"
$ws.Range("M497").Value = "N/A"
$ws.Range("N497").Value = "N/A"
$ws.Range("O497").NumberFormat = "@"
$ws.Range("O497").Value = "2025"

# --- Row 498 ---
$ws.Range("A498").Value = "The Thal from G.R.A.C.E"
$ws.Range("E498").Value = "N/A"
$ws.Range("F498").NumberFormat = "@"
$ws.Range("F498").Value = "2"
$ws.Range("G498").Value = "N/A"
$ws.Range("H498").Value = 84
$ws.Range("I498").Value = "N, /, A"
$ws.Range("J498").Value = "N/A"
$ws.Range("K498").Value = "N, /, A"
$ws.Range("L498").Value = "N/A"
$ws.Range("M498").Value = "N/A"
$ws.Range("N498").Value = "N/A"
$ws.Range("O498").NumberFormat = "@"
$ws.Range("O498").Value = "2025"

# --- Row 499 ---
$ws.Range("A499").Value = "Allegiance"
$ws.Range("E499").Value = "N/A"
$ws.Range("F499").NumberFormat = "@"
$ws.Range("F499").Value = "3"
$ws.Range("G499").Value = "N/A"
$ws.Range("H499").Value = 86
$ws.Range("I499").Value = "N, /, A"
$ws.Range("J499").Value = "N/A"
$ws.Range("K499").Value = "N, /, A"
$ws.Range("L499").Value = "N/A"
$ws.Range("M499").Value = "N/A"
$ws.Range("N499").Value = "N/A"
$ws.Range("O499").NumberFormat = "@"
$ws.Range("O499").Value = "2025"

# --- Row 500 ---
$ws.Range("A500").Value = "Spoil of War"
$ws.Range("E500").Value = "Pursuit"
$ws.Range("F500").NumberFormat = "@"
$ws.Range("F500").Value = "1"
$ws.Range("G500").NumberFormat = "@"
$ws.Range("G500").Value = "1"
$ws.Range("H500").Value = 76
$ws.Range("I500").Value = "Alex, Cass"
$ws.Range("K500").Value = "8th Dr, Hieronyma Friend"
$ws.Range("L500").Value = "Sontarans"
$ws.Range("M500").Value = "Mark Wright"
$ws.Range("N500").Value = "Ken Bentley"
$ws.Range("O500").NumberFormat = "@"
$ws.Range("O500").Value = "2025"

# --- Row 501 ---
$ws.Range("A501").Value = "The Tale of Alex"
$ws.Range("E501").Value = "Pursuit"
$ws.Range("F501").NumberFormat = "@"
$ws.Range("F501").Value = "2"
$ws.Range("G501").NumberFormat = "@"
$ws.Range("G501").Value = "1"
$ws.Range("H501").Value = 76
$ws.Range("I501").Value = "Eighth Doctor, Alex, Cass, Hieronyma Friend"
$ws.Range("M501").Value = "Katharine Armitage"
$ws.Range("N501").Value = "Ken Bentley"
$ws.Range("O501").NumberFormat = "@"
$ws.Range("O501").Value = "2025"

# --- Row 502 ---
$ws.Range("A502").Value = "See-Saw"
$ws.Range("E502").Value = "Pursuit"
$ws.Range("F502").NumberFormat = "@"
$ws.Range("F502").Value = "3"
$ws.Range("G502").NumberFormat = "@"
$ws.Range("G502").Value = "1"
$ws.Range("H502").Value = 67
$ws.Range("I502").Value = "Eighth Doctor, Alex, Cass, Hieronyma Friend"
$ws.Range("M502").Value = "James Moran"
$ws.Range("N502").Value = "Ken Bentley"
$ws.Range("O502").NumberFormat = "@"
$ws.Range("O502").Value = "2025"

# --- Row 503 ---
$ws.Range("A503").Value = "The First Forest"
# B503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
# C503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
# D503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
$ws.Range("E503").Value = "Pursuit"
$ws.Range("F503").NumberFormat = "@"
$ws.Range("F503").Value = "4"
$ws.Range("G503").NumberFormat = "@"
$ws.Range("G503").Value = "1"
$ws.Range("H503").Value = 76
$ws.Range("I503").Value = "Eighth Doctor, Alex, Cass, Hieronyma Friend"
# J503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
# K503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
# L503: empty placeholder cell in source data (left unset; COM Value="" cannot create a present-but-empty cell)
$ws.Range("M503").Value = "Tim Foley"
$ws.Range("N503").Value = "Ken Bentley"
$ws.Range("O503").NumberFormat = "@"
$ws.Range("O503").Value = "2025"
